$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.710.37"
$ws.Range("E2").Value = "  -0.13%  "
$ws.Range("D3").Value = "1.914.77"
$ws.Range("E3").Value = "  +0.80%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  +0.23%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "239.74"
$ws.Range("E5").Value = "  -2.30%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.002"
$ws.Range("E6").Value = "  +0.27%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4920"
$ws.Range("E7").Value = "  -0.23%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2975"
$ws.Range("E8").Value = "  +0.60%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06767"
$ws.Range("E9").Value = "  -0.42%  "
$ws.Range("D10").Value = "1.891.62"
$ws.Range("E10").Value = "  -0.47%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "17.17"
$ws.Range("E11").Value = "  -0.76%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07369"
$ws.Range("E12").Value = "  +1.51%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.174"
$ws.Range("E13").Value = "  +2.40%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "88.72"
$ws.Range("E14").Value = "  -2.70%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6714"
$ws.Range("E15").Value = "  -1.45%  "
$ws.Range("D16").Value = "30.685.15"
$ws.Range("E16").Value = "  -0.16%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000007946"
$ws.Range("E17").Value = "  -1.08%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.53"
$ws.Range("E18").Value = "  +2.49%  "
$ws.Range("E19").Value = "  +0.31%  "
$ws.Range("D20").Value = "2.170.23"
$ws.Range("E20").Value = "  +1.59%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.407"
$ws.Range("E21").Value = "  +12.18%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.002"
$ws.Range("E22").Value = "  +0.27%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "200.97"
$ws.Range("E23").Value = "  +3.88%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.298"
$ws.Range("E24").Value = "  +2.43%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.663"
$ws.Range("E25").Value = "  +2.48%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "164.46"
$ws.Range("E26").Value = "  +5.31%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.72"
$ws.Range("E27").Value = "  -2.52%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.956"
$ws.Range("E28").Value = "  +2.14%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.486"
$ws.Range("E29").Value = "  +6.17%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.359"
$ws.Range("E30").Value = "  +0.44%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09141"
$ws.Range("E31").Value = "  +0.32%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.061"
$ws.Range("E32").Value = "  +0.91%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05296"
$ws.Range("E33").Value = "  +1.03%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7392"
$ws.Range("E34").Value = "  -1.21%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.114"
$ws.Range("E35").Value = "  -0.09%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.726"
$ws.Range("E36").Value = "  -1.26%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01832"
$ws.Range("E37").Value = "  -0.85%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.715"
$ws.Range("E38").Value = "  +0.61%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.9223"
$ws.Range("E39").Value = "  -1.98%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.082"
$ws.Range("E40").Value = "  -2.93%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "74.87"
$ws.Range("E41").Value = "  +29.46%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.4465"
$ws.Range("E42").Value = "  +0.83%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.002"
$ws.Range("E45").Value = "  +0.25%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.1390"
$ws.Range("E46").Value = "  +2.65%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.639"
$ws.Range("E47").Value = "  +0.01%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "35.55"
$ws.Range("E48").Value = "  +5.35%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.096"
$ws.Range("E49").Value = "  +3.52%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05880"
$ws.Range("E50").Value = "  +0.22%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4021"
$ws.Range("E51").Value = "  +1.60%  "

# Row 43 becomes Quant (previously row 44 contents with updated price/change)
$ws.Range("B43").Value = "Quant"
$ws.Range("C43").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "107.13"
$ws.Range("E43").Value = "  +1.50%  "

# Row 44 becomes FraxShare (previously row 43 contents with updated price/change)
$ws.Range("B44").Value = "FraxShare"
$ws.Range("C44").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.943"
$ws.Range("E44").Value = "  +2.90%  "
